$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# New log entry added for 31.10.2018 (row 24 of the time-tracking table)
$ws.Range("A24").Value2 = 43404
$ws.Range("B24").Value2 = 0.48402777777777778
$ws.Range("C24").Value2 = 0.80555555555555547
$ws.Range("D24").Formula = "=C24-B24"
$ws.Range("E24").Value2 = 3
$ws.Range("F24").Value2 = "4 h midPointin ja testipalvelimen liittämistä niin että midPointista saadaan lisättyä käyttäjiä palvelimeen oikeilla asetuksilla (mm. admin käyttäjät saavat sudo oikeudet ja normaalikäyttäjät eivät pysty mm. muuttamaan asetustiedostoja jne.) Testattiin liittää useita erilaisia käyttäjiä palvelimeen midPointin kautta sekä ryhmä ja rooli jakoa käyttäjille. 1h 30 min Otin selvää midPointin lokeista /var/log/authlog, /var/log/syslog --> ei löytynyt mm. käyttäjien lisäämis tietoja ja ajankohtia. Katsoin myös midPoint GUI:sta lokeja jos löytyisi kiinnostavampaa loki tietoa, oli hieman epäselvää joten katsoin midPointin omasta dokumentaatiosta apua lokien selaamiseen --> https://wiki.evolveum.com/display/midPoint/Log+Viewer"

$ws.Rows.Item(24).RowHeight = 195

# Reflect the selection left behind after entering the row
[void]$ws.Range("C24").Select()
